# Update countries & provincias Spain
# Applies the COVID-19 "Pais" data refresh described by the commit diff:
#  - Updated timestamp banner in A1
#  - Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados,
#    Casos criticos, Muertes hoy, Muertes) for several countries
#  - Re-sorted two pairs of countries whose case counts changed enough to swap
#    their ranking position in the (descending, by total cases) table:
#      * Cuba overtook Afganistan (rows 82/83)
#      * Birmania overtook Aruba & Guayana Francesa (rows 140/141/142)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "last refreshed" banner
$ws.Range("A1").Value = "Datos actualizados a 18 de Abril de 2020 a las 17:52"

# Row -> country name shown in column A (only rows whose country changed are listed)
$countries = @{
    82  = "Cuba"
    83  = "Afganistan"
    140 = "Birmania"
    141 = "Aruba"
    142 = "Guayana Francesa"
}

foreach ($row in $countries.Keys) {
    $ws.Cells.Item([int]$row, 1).Value = $countries[$row]
}

# Row -> statistics in columns B..H
#        B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#        F=Casos criticos, G=Muertes hoy, H=Muertes
$stats = @{
    4   = @(712719, 2984, 63778, 611652, 13544, 135, 37289)
    18  = @(27404,  326,  16400, 9638,   386,   39,  1366)
    28  = @(9730,   478,  4035,  5569,   360,   10,  126)
    29  = @(8742,   363,  981,   7414,   160,   15,  347)
    51  = @(3537,   57,   601,   2864,   32,    0,   72)
    82  = @(986,    63,   227,   727,    17,    1,   32)
    83  = @(933,    27,   112,   791,    0,     0,   30)
    140 = @(98,     10,   5,     88,     0,     1,   5)
    141 = @(96,     0,    43,    51,     1,     0,   2)
    142 = @(96,     0,    64,    32,     2,     0,   0)
}

foreach ($row in $stats.Keys) {
    $values = $stats[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item([int]$row, $i + 2).Value = $values[$i]
    }
}
